$wb = $excel.ActiveWorkbook

# Add the new "OOTB Domain Groups" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "OOTB Domain Groups"

# Populate the OOTB Domain Groups data (provider name + associated domains per row)
$row1 = @('Gmail', 'gmail.com', 'googlemail.com', 'googlemail.co.uk')
for ($i = 0; $i -lt $row1.Length; $i++) { $ws3.Cells.Item(1, $i + 1).Value = $row1[$i] }
$row2 = @('Microsoft', 'live.com', 'msn.com', 'hotmail.ca', 'hotmail.com', 'hotmail.de', 'hotmail.dk', 'hotmail.co.jp', 'hotmail.it', 'hotmail.es', 'hotmail.fr', 'hotmail.co.uk', 'hotmail.co.kr', 'hotmail.com.au', 'hotmail.com.ar', 'hotmail.co.il', 'hotmail.com.br', 'hotmail.com.tr', 'hotmail.co.th', 'hotmail.jp', 'hotmail.se', 'live.at', 'live.be', 'live.ca', 'live.cl', 'live.cn', 'live.co.kr', 'live.co.uk', 'live.com.ar', 'live.com.au', 'live.com.mx', 'live.com.my', 'live.com.sg', 'live.de', 'live.dk', 'live.fr', 'live.hk', 'live.ie', 'live.in', 'live.it', 'live.jp', 'live.nl', 'live.no', 'live.ru', 'live.se', 'outlook.com', 'live.co.uk', 'hotmail.gr', 'windowslive.com', 'xbox.com', 'hotmail.cl', 'live.at', 'live.jp', 'live.ca', 'hotmail.ca', 'hotmail.ch', 'live.fr', 'live.it', 'live.nl', 'outlook.ie', 'outlook.com.br', 'live.com.pt', 'live.be', 'live.co.za', 'mts.net')
for ($i = 0; $i -lt $row2.Length; $i++) { $ws3.Cells.Item(2, $i + 1).Value = $row2[$i] }
$row3 = @('Yahoo', 'yahoo.com', 'rocketmail.com', 'rogers.com', 'sky.com', 'talk21.com', 'y7mail.com', 'yahoo.at', 'yahoo.be', 'yahoo.bg', 'yahoo.ca', 'yahoo.cl', 'yahoo.co.hu', 'yahoo.co.id', 'yahoo.co.il', 'yahoo.co.in', 'yahoo.co.jp', 'yahoo.co.kr', 'yahoo.com.ar', 'yahoo.com.au', 'yahoo.com.biz', 'yahoo.com.br', 'yahoo.com.cn', 'yahoo.com.co', 'yahoo.com.hk', 'yahoo.com.hr', 'yahoo.com.in', 'yahoo.com.jp', 'yahoo.com.kr', 'yahoo.com.mx', 'yahoo.com.my', 'yahoo.com.net', 'yahoo.com.pe', 'yahoo.com.ph', 'yahoo.com.sg', 'yahoo.com.tr', 'yahoo.com.tw', 'yahoo.com.ua', 'yahoo.com.ve', 'yahoo.com.vn', 'yahoo.co.nz', 'yahoo.co.th', 'yahoo.co.uk', 'yahoo.co.za', 'yahoo.cz', 'yahoo.de', 'yahoo.dk', 'yahoo.ee', 'yahoo.es', 'yahoo.fi', 'yahoo.fr', 'yahoogroups.co.kr', 'yahoogroups.com.cn', 'yahoogroups.com.sg', 'yahoogroups.com.tw', 'yahoogrupper.dk', 'yahoogruppi.it', 'yahoo.gr', 'yahoo.hr', 'yahoo.hu', 'yahoo.ie', 'yahoo.in', 'yahoo.it', 'yahoo.lt', 'yahoo.lv', 'yahoo.nl', 'yahoo.no', 'yahoo.pl', 'yahoo.pt', 'yahoo.ro', 'yahoo.rs', 'yahoo.se', 'yahoo.si', 'yahoo.sk', 'yahooxtra.co.nz', 'ymail.com', 'aol.com', 'aim.com', 'compuserve.com', 'cs.com', 'netscape.com', 'netscape.net', 'wmconnect.com', 'aol.co.uk', 'aol.in', 'aol.de', 'aol.fr', 'aol.nl', 'aol.pl', 'aol.jp', 'aol.es', 'aol.it', 'aol.com.ar', 'aol.fi', 'aol.cl', 'aol.com.co', 'aol.com.ve', 'aol.com.au', 'aol.at', 'aol.be', 'aol.com.br', 'aol.cz', 'aol.dk', 'myaol.jp', 'aolnorge.no', 'aolpolska.pl', 'aolpolcka.pl', 'aolpoland.pl', 'aol.ru', 'aol.kr', 'aol.se', 'aol.ch', 'aol.com.tr', 'aol.co.nz', 'aolchina.com', 'aol.hk', 'aol.tw', 'luckymail.com', 'verizon.net', 'aol.com.mx', 'bellatlantic.net', 'citlink.net', 'frontier.com', 'frontiernet.net', 'games.com', 'goowy.com', 'gte.net', 'love.com', 'verizon.net.in', 'wild4music.com', 'wow.com', 'yahoo.cn', 'yahoo.ne.jp', 'yahoogroups.ca', 'yahoogroups.co.in', 'yahoogroups.co.uk', 'yahoogroups.com', 'yahoogroups.com.au', 'yahoogroups.com.hk', 'yahoogroups.de', 'ybb.ne.jp', 'ygm.com')
for ($i = 0; $i -lt $row3.Length; $i++) { $ws3.Cells.Item(3, $i + 1).Value = $row3[$i] }
$row4 = @('Apple', 'mac.com', 'icloud.com', 'apple.com', 'me.com')
for ($i = 0; $i -lt $row4.Length; $i++) { $ws3.Cells.Item(4, $i + 1).Value = $row4[$i] }
$row5 = @('Comcast', 'comcast.net')
for ($i = 0; $i -lt $row5.Length; $i++) { $ws3.Cells.Item(5, $i + 1).Value = $row5[$i] }
$row6 = @('Orange', 'orange.fr', 'orange.com', 'wanadoo.fr', 'francetelecom.com', 'voila.fr', 'voila.com')
for ($i = 0; $i -lt $row6.Length; $i++) { $ws3.Cells.Item(6, $i + 1).Value = $row6[$i] }
$row7 = @('La Poste', 'laposte.net')
for ($i = 0; $i -lt $row7.Length; $i++) { $ws3.Cells.Item(7, $i + 1).Value = $row7[$i] }
$row8 = @('Italia Online', 'libero.it', 'inwind.it', 'iol.it', 'blu.it', 'giallo.it', 'virgilio.it')
for ($i = 0; $i -lt $row8.Length; $i++) { $ws3.Cells.Item(8, $i + 1).Value = $row8[$i] }
$row9 = @('WP', 'wp.pl', 'o2.pl')
for ($i = 0; $i -lt $row9.Length; $i++) { $ws3.Cells.Item(9, $i + 1).Value = $row9[$i] }
$row10 = @('United Internet', 'web.de', 'gmx.de', 'gmx.ch', 'gmx.net', 'gmx.com', 'gmx.at', 'gmx.fr', 'mail.com', '1and1.com', '1und1.de')
for ($i = 0; $i -lt $row10.Length; $i++) { $ws3.Cells.Item(10, $i + 1).Value = $row10[$i] }
$row11 = @('Bigpond', 'bigpond.com', 'bigpond.net.au', 'bigpond.com.au', 'telstra.com', 'bigpond.net')
for ($i = 0; $i -lt $row11.Length; $i++) { $ws3.Cells.Item(11, $i + 1).Value = $row11[$i] }
$row12 = @('Docomo', 'docomo.ne.jp')
for ($i = 0; $i -lt $row12.Length; $i++) { $ws3.Cells.Item(12, $i + 1).Value = $row12[$i] }
$row13 = @('Softbank', 'softbank.ne.jp', 'c.vodafone.ne.jp', 'd.vodafone.ne.jp', 'h.vodafone.ne.jp', 'k.vodafone.ne.jp', 'n.vodafone.ne.jp', 'q.vodafone.ne.jp', 'r.vodafone.ne.jp', 's.vodafone.ne.jp', 't.vodafone.ne.jp', 'jp-c.ne.jp', 'jp-d.ne.jp', 'jp-h.ne.jp', 'jp-k.ne.jp', 'jp-n.ne.jp', 'jp-q.ne.jp', 'jp-r.ne.jp', 'jp-s.ne.jp', 'jp-t.ne.jp')
for ($i = 0; $i -lt $row13.Length; $i++) { $ws3.Cells.Item(13, $i + 1).Value = $row13[$i] }
$row14 = @('KDDI', 'au.com', 'ezweb.ne.jp', 'uqmobile.jp')
for ($i = 0; $i -lt $row14.Length; $i++) { $ws3.Cells.Item(14, $i + 1).Value = $row14[$i] }

# Default selection on the new sheet
$null = $ws3.Range("E8").Select()

# Fix up the font/style of A1 on "Custom Domain Group" (Menlo -> Calibri, matches B1 style)
$ws2 = $wb.Worksheets.Item("Custom Domain Group")
$ws2.Range("A1").Font.Name = $ws2.Range("B1").Font.Name
$ws2.Range("A1").Font.Color = $ws2.Range("B1").Font.Color
$ws2.Range("A1").Font.Size = $ws2.Range("B1").Font.Size

# Update the active selection on "Custom Domain Group"
$null = $ws2.Range("C6").Select()

# Make "Warmup Plan" the active/selected sheet (tab selection)
$ws1 = $wb.Worksheets.Item("Warmup Plan")
$null = $ws1.Activate()
